# feat: add 2022-Q3 data
#
# - Insert a brand-new worksheet "2022-Q3" right after "总计" (i.e. before
#   the existing "2022-Q2" sheet), carrying the new quarter's fund snapshot.
# - Update the "总计" (totals) sheet: every existing row shifts down by one
#   (the new quarter becomes the latest entry) and a new last row is added
#   for the oldest quarter that fell off the bottom.
# - The other existing quarter sheets ("2022-Q2", "2021-Q3", "2021-Q2") are
#   left untouched; they simply shift tab position because of the insert.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) sheet.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 4.06

$total.Range("B3").Value = "2022-Q2"
$total.Range("D3").Value = 5.32

$total.Range("B4").Value = "2021-Q3"
$total.Range("D4").Value = 2.8

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 1.88

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" sheet before the current "2022-Q2" sheet.
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "164906"
$q3.Range("C2").Value = "交银施罗德中证海外中国互联网指数（QDII-LOF）"
$q3.Range("D2").Value = "97.68"
$q3.Range("E2").Value = "91.19"
$q3.Range("F2").Value = "4.16"
$q3.Range("G2").Value = "4.0635"
$q3.Range("H2").Value = 9
